# Generate Report for Handback
# Adds a new handback row (file 6edb2551-9639-414c-aa0a-3a39fd44bb14.md) to
# each of the three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newMd   = "6edb2551-9639-414c-aa0a-3a39fd44bb14.md"
$newPath = "e2e\6edb2551-9639-414c-aa0a-3a39fd44bb14.md"

$zhXlf = "6edb2551-9639-414c-aa0a-3a39fd44bb14.dc0c91972cb0f44849542d8152c9644b80388c93.zh-cn.xlf"
$deXlf = "6edb2551-9639-414c-aa0a-3a39fd44bb14.dc0c91972cb0f44849542d8152c9644b80388c93.de-de.xlf"

# -----------------------------------------------------------------
# Sheet "Overview"
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $newMd
$wsOverview.Range("B4").Value = $newPath
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc0c91972cb0f44849542d8152c9644b80388c93/e2e/6edb2551-9639-414c-aa0a-3a39fd44bb14.md", "", "", $newPath) | Out-Null
$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-17 22:43:50"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -----------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $newMd
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/dc0c91972cb0f44849542d8152c9644b80388c93/e2e/6edb2551-9639-414c-aa0a-3a39fd44bb14.md", "", "", $newMd) | Out-Null
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = "2016-08-17 22:43:45"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $newMd
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/dc0c91972cb0f44849542d8152c9644b80388c93/e2e/6edb2551-9639-414c-aa0a-3a39fd44bb14.md", "", "", $newMd) | Out-Null
$wsZh.Range("I4").Style = "HyperLink"
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = "2016-08-17 22:44:07"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

# -----------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $newMd
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dc0c91972cb0f44849542d8152c9644b80388c93/e2e/6edb2551-9639-414c-aa0a-3a39fd44bb14.md", "", "", $newMd) | Out-Null
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = "2016-08-17 22:43:50"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $newMd
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dc0c91972cb0f44849542d8152c9644b80388c93/e2e/6edb2551-9639-414c-aa0a-3a39fd44bb14.md", "", "", $newMd) | Out-Null
$wsDe.Range("I4").Style = "HyperLink"
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = "2016-08-17 22:44:15"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

# -----------------------------------------------------------------
# Extend the tables (autofilter + ref) to include the new row.
# -----------------------------------------------------------------
$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G4"))
$wsZh.ListObjects.Item(1).Resize($wsZh.Range("A1:P4"))
$wsDe.ListObjects.Item(1).Resize($wsDe.Range("A1:P4"))
